$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark the first task ("abgeschlossen" = "completed") in the Bemerkung column
$ws.Range("F2").Value = "abgeschlossen"

# Move the active selection to D2, matching the saved workbook state
$null = $ws.Range("D2").Select()
